$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$meta = $wb.Worksheets.Item("Metadata")

# Version: 2.2.0-ballot -> 2.1.0
$meta.Range("B3").Value = "2.1.0"

# Date: 2025-12-19T08:32:44+00:00 -> 2025-12-19T08:44:55+00:00
$meta.Range("B8").Value = "2025-12-19T08:44:55+00:00"

# Base Definition: drop the "|4.0.1" version suffix
$meta.Range("B18").Value = "http://hl7.org/fhir/StructureDefinition/Extension"

# --- Elements sheet updates ---
$elements = $wb.Worksheets.Item("Elements")

# Binding Value Set: drop the "|20250624152100" version suffix
$elements.Range("Z6").Value = "https://smt.esante.gouv.fr/fhir/ValueSet/jdv-motif-non-realisation-evenement-cisis"

# Column Z width shrinks now that the long value-set URI lost its version suffix
# (target stored width 65.45703125 character-units; feed the equivalent
# ColumnWidth so the engine's internal pixel-rounding lands as close as possible)
$elements.Columns.Item(26).ColumnWidth = (65.45703125 - 5/6)
